$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 15-22 (old rows beyond new data range)
$ws.Range("A15:G22").EntireRow.Delete() | Out-Null

# Update rows 2-14 with new year labels and data (years 2010-2022)
$ws.Range("A2").Value = "2010年"
$ws.Range("B2").Value = 5.0459447447303
$ws.Range("C2").Value = 47.4427032280895
$ws.Range("D2").Value = -1.14988538601278
$ws.Range("E2").Value = -10.8113889221427
$ws.Range("F2").Value = 6.73981170552388
$ws.Range("G2").Value = 63.3686856940531

$ws.Range("A3").Value = "2011年"
$ws.Range("B3").Value = 6.27423040605307
$ws.Range("C3").Value = 65.693023272134
$ws.Range("D3").Value = -0.6474544919201149
$ws.Range("E3").Value = -6.77903746797724
$ws.Range("F3").Value = 3.9240562648977
$ws.Range("G3").Value = 41.0860141958433

$ws.Range("A4").Value = "2012年"
$ws.Range("B4").Value = 4.35620880765038
$ws.Range("C4").Value = 55.3961699506632
$ws.Range("D4").Value = 0.194198172006589
$ws.Range("E4").Value = 2.46954069825399
$ws.Range("F4").Value = 3.31332946914463
$ws.Range("G4").Value = 42.1342893510828

$ws.Range("A5").Value = "2013年"
$ws.Range("B5").Value = 3.89560592278008
$ws.Range("C5").Value = 50.1613524562137
$ws.Range("D5").Value = -0.253045184019436
$ws.Range("E5").Value = -3.25830921159706
$ws.Range("F5").Value = 4.12358935892779
$ws.Range("G5").Value = 53.0969567553834

$ws.Range("A6").Value = "2014年"
$ws.Range("B6").Value = 4.18189662923193
$ws.Range("C6").Value = 56.316048057266
$ws.Range("D6").Value = -0.094197061834455
$ws.Range("E6").Value = -1.26851683134423
$ws.Range("F6").Value = 3.33806408878006
$ws.Range("G6").Value = 44.9524687740783

$ws.Range("A7").Value = "2015年"
$ws.Range("B7").Value = 4.85707638314603
$ws.Range("C7").Value = 68.97954160093239
$ws.Range("D7").Value = 0.593241593878498
$ws.Range("E7").Value = 8.42513684700169
$ws.Range("F7").Value = 1.59101090168226
$ws.Range("G7").Value = 22.5953215520656

$ws.Range("A8").Value = "2016年"
$ws.Range("B8").Value = 4.51806962864869
$ws.Range("C8").Value = 65.969141480161
$ws.Range("D8").Value = -0.799321371989767
$ws.Range("E8").Value = -11.6710340944171
$ws.Range("F8").Value = 3.13001394832964
$ws.Range("G8").Value = 45.7018926142563

$ws.Range("A9").Value = "2017年"
$ws.Range("B9").Value = 3.88124001078923
$ws.Range("C9").Value = 55.8676814771848
$ws.Range("D9").Value = 0.323844617226383
$ws.Range("E9").Value = 4.66151226747386
$ws.Range("F9").Value = 2.74211616529859
$ws.Range("G9").Value = 39.4708062553414

$ws.Range("A10").Value = "2018年"
$ws.Range("B10").Value = 4.31883411587706
$ws.Range("C10").Value = 63.984871538296
$ws.Range("D10").Value = -0.485021427698422
$ws.Range("E10").Value = -7.18574340017275
$ws.Range("F10").Value = 2.91596114436263
$ws.Range("G10").Value = 43.2008718618766

$ws.Range("A11").Value = "2019年"
$ws.Range("B11").Value = 3.48476514656898
$ws.Range("C11").Value = 58.5625533751102
$ws.Range("D11").Value = 0.7469647525442
$ws.Range("E11").Value = 12.5529731130562
$ws.Range("F11").Value = 1.71877079754881
$ws.Range("G11").Value = 28.8844735118336

$ws.Range("A12").Value = "2020年"
$ws.Range("B12").Value = -0.153173366486956
$ws.Range("C12").Value = -6.84224011054964
$ws.Range("D12").Value = 0.567040796794255
$ws.Range("E12").Value = 25.3296599345427
$ws.Range("F12").Value = 1.82477611350502
$ws.Range("G12").Value = 81.51258017600711

$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 4.92308804287977
$ws.Range("C13").Value = 58.261396957157
$ws.Range("D13").Value = 1.85263953694034
$ws.Range("E13").Value = 21.9247282478147
$ws.Range("F13").Value = 1.67427242017988
$ws.Range("G13").Value = 19.8138747950282

$ws.Range("A14").Value = "2022年"
$ws.Range("B14").Value = 0.981413118413281
$ws.Range("C14").Value = 32.8231812178355
$ws.Range("D14").Value = 0.511549698598217
$ws.Range("E14").Value = 17.1086855718467
$ws.Range("F14").Value = 1.49703718298851
$ws.Range("G14").Value = 50.0681332103182

